$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<zero>"
$ws.Range("C2").Value = 19

# Row 3
$ws.Range("C3").Value = 25

# Row 4
$ws.Range("C4").Value = 24

# Row 5
$ws.Range("B5").Value = "<lone>"
$ws.Range("C5").Value = 37

# Row 6
$ws.Range("C6").Value = 30

# Row 7
$ws.Range("B7").Value = "<up>"
$ws.Range("C7").Value = 29

# Row 8
$ws.Range("C8").Value = 32

# Row 9
$ws.Range("B9").Value = "<now>"
$ws.Range("C9").Value = 19

# Row 10
$ws.Range("C10").Value = 28

# Row 12
$ws.Range("C12").Value = 26

# Row 13
$ws.Range("B13").Value = "<mace>"
$ws.Range("C13").Value = 28

# Row 14
$ws.Range("C14").Value = 32

# Row 15
$ws.Range("B15").Value = "<can>"
$ws.Range("C15").Value = 6
